$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scraped values updated by the GitHub Action).
# Price (D) and Volume (E) cells are forced to Text format right before each write so
# numeric-looking values (e.g. "2.30", "1.00") keep their original text appearance
# instead of being auto-converted to numbers and losing significant trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.359.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.713.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.38%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.712.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.329.71"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.712.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.412.88"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.74"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +8.60%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.14"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000141"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.19"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.43"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.850.44"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.650.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.22%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "48.78"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "426.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.92%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.43"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.57"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.745.47"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.11%  "
